$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 11500
$ws.Range("J29").Value = 12550
$ws.Range("L29").Value = 37650
$ws.Range("N29").Value = -38212
$ws.Range("H38").Value = 3792.7144
$ws.Range("I38").Value = 2758.25
$ws.Range("J38").Value = 9999.5
$ws.Range("K38").Value = 8274.75
$ws.Range("L38").Value = 29998.5
$ws.Range("M38").Value = -7902.75
$ws.Range("N38").Value = -30742.5
$ws.Range("H40").Value = 9465.166999999999
$ws.Range("I40").Value = 11749.75
$ws.Range("J40").Value = 4896
$ws.Range("K40").Value = 11749.75
$ws.Range("L40").Value = 4896
$ws.Range("M40").Value = -11574.75
$ws.Range("N40").Value = -5246
$ws.Range("H80").Value = 2595.762
$ws.Range("I80").Value = 976.75
$ws.Range("J80").Value = 3592.077
$ws.Range("K80").Value = 2930.25
$ws.Range("L80").Value = 10776.231
$ws.Range("M80").Value = -1932.25
$ws.Range("N80").Value = -12772.231
$ws.Range("H83").Value = 2595.762
$ws.Range("I83").Value = 976.75
$ws.Range("J83").Value = 3592.077
$ws.Range("K83").Value = 8790.75
$ws.Range("L83").Value = 32328.693
$ws.Range("M83").Value = -3798.75
$ws.Range("N83").Value = -42312.693
$ws.Range("H121").Value = 4239.5
$ws.Range("J121").Value = 4239.5
$ws.Range("L121").Value = 12718.5
$ws.Range("N121").Value = -16212.5
$ws.Range("H123").Value = 58200
$ws.Range("J123").Value = 58200
$ws.Range("L123").Value = 58200
$ws.Range("N123").Value = -68000
$ws.Range("H131").Value = 3381.2563
$ws.Range("I131").Value = 1580.9615
$ws.Range("J131").Value = 6981.846
$ws.Range("K131").Value = 4742.8845
$ws.Range("L131").Value = 20945.538
$ws.Range("M131").Value = 297.1154999999999
$ws.Range("N131").Value = -31025.538
$ws.Range("H137").Value = 10802.923
$ws.Range("I137").Value = 1678.4286
$ws.Range("J137").Value = 21448.166
$ws.Range("K137").Value = 5035.2858
$ws.Range("L137").Value = 64344.49800000001
$ws.Range("M137").Value = -2485.2858
$ws.Range("N137").Value = -69444.49800000001
$ws.Range("H138").Value = 11093.204
$ws.Range("I138").Value = 3128
$ws.Range("J138").Value = 12885.375
$ws.Range("K138").Value = 9384
$ws.Range("L138").Value = 38656.125
$ws.Range("M138").Value = -4244
$ws.Range("N138").Value = -48936.125
$ws.Range("H141").Value = 5036.4287
$ws.Range("I141").Value = 5481.8184
$ws.Range("J141").Value = 3403.3333
$ws.Range("K141").Value = 16445.4552
$ws.Range("L141").Value = 10209.9999
$ws.Range("M141").Value = -11265.4552
$ws.Range("N141").Value = -20569.9999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3730.9167
$ws.Range("I2").Value = 2584.6667
$ws.Range("J2").Value = 7169.6665
$ws.Range("K2").Value = 2584.6667
$ws.Range("L2").Value = 7169.6665
$ws.Range("M2").Value = -2471.6667
$ws.Range("N2").Value = -7395.6665
$ws.Range("H32").Value = 14928836
$ws.Range("I32").Value = 16950726
$ws.Range("K32").Value = 16950726
$ws.Range("M32").Value = -16950439
$ws.Range("H61").Value = 62503490
$ws.Range("I61").Value = 100003224
$ws.Range("J61").Value = 3916.1667
$ws.Range("K61").Value = 100003224
$ws.Range("L61").Value = 3916.1667
$ws.Range("M61").Value = -100003012
$ws.Range("N61").Value = -4340.1667
$ws.Range("H74").Value = 37080348
$ws.Range("I74").Value = 45506700
$ws.Range("J74").Value = 4399.6
$ws.Range("K74").Value = 45506700
$ws.Range("L74").Value = 4399.6
$ws.Range("M74").Value = -45505826
$ws.Range("N74").Value = -6147.6
$ws.Range("H77").Value = 37080348
$ws.Range("I77").Value = 45506700
$ws.Range("J77").Value = 4399.6
$ws.Range("K77").Value = 227533500
$ws.Range("L77").Value = 21998
$ws.Range("M77").Value = -227529132
$ws.Range("N77").Value = -30734
$ws.Range("H116").Value = 3730.9167
$ws.Range("I116").Value = 2584.6667
$ws.Range("J116").Value = 7169.6665
$ws.Range("K116").Value = 2584.6667
$ws.Range("L116").Value = 7169.6665
$ws.Range("M116").Value = -290.6667000000002
$ws.Range("N116").Value = -11757.6665
$ws.Range("H122").Value = 10755521
$ws.Range("I122").Value = 2356.5454
$ws.Range("J122").Value = 37041030
$ws.Range("K122").Value = 7069.6362
$ws.Range("L122").Value = 111123090
$ws.Range("M122").Value = -4619.6362
$ws.Range("N122").Value = -111127990
$ws.Range("H132").Value = 32307330
$ws.Range("I132").Value = 4662.913
$ws.Range("J132").Value = 125177496
$ws.Range("K132").Value = 13988.739
$ws.Range("L132").Value = 375532488
$ws.Range("M132").Value = -11458.739
$ws.Range("N132").Value = -375537548
$ws.Range("H136").Value = 62503490
$ws.Range("I136").Value = 100003224
$ws.Range("J136").Value = 3916.1667
$ws.Range("K136").Value = 300009672
$ws.Range("L136").Value = 11748.5001
$ws.Range("M136").Value = -300007122
$ws.Range("N136").Value = -16848.5001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3730.9167
$ws.Range("I3").Value = 2584.6667
$ws.Range("J3").Value = 7169.6665
$ws.Range("K3").Value = 2584.6667
$ws.Range("L3").Value = 7169.6665
$ws.Range("M3").Value = -2470.6667
$ws.Range("N3").Value = -7397.6665
$ws.Range("H94").Value = 1567.3043
$ws.Range("I94").Value = 902.5
$ws.Range("K94").Value = 902.5
$ws.Range("M94").Value = -451.5
$ws.Range("H99").Value = 4229.7144
$ws.Range("I99").Value = 3222
$ws.Range("J99").Value = 6749
$ws.Range("K99").Value = 3222
$ws.Range("L99").Value = 6749
$ws.Range("M99").Value = -1724
$ws.Range("N99").Value = -9745
$ws.Range("H134").Value = 3166.1052
$ws.Range("I134").Value = 2991.6
$ws.Range("J134").Value = 3820.5
$ws.Range("K134").Value = 8974.799999999999
$ws.Range("L134").Value = 11461.5
$ws.Range("M134").Value = -6439.799999999999
$ws.Range("N134").Value = -16531.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 994.4167
$ws.Range("I16").Value = 1028.3334
$ws.Range("J16").Value = 892.6667
$ws.Range("K16").Value = 1028.3334
$ws.Range("L16").Value = 892.6667
$ws.Range("M16").Value = -741.3334
$ws.Range("N16").Value = -1466.6667
$ws.Range("H31").Value = 35716996
$ws.Range("I31").Value = 2330.8
$ws.Range("J31").Value = 62502996
$ws.Range("K31").Value = 2330.8
$ws.Range("L31").Value = 62502996
$ws.Range("M31").Value = -2035.8
$ws.Range("N31").Value = -62503586
$ws.Range("H34").Value = 35716996
$ws.Range("I34").Value = 2330.8
$ws.Range("J34").Value = 62502996
$ws.Range("K34").Value = 2330.8
$ws.Range("L34").Value = 62502996
$ws.Range("M34").Value = -2128.8
$ws.Range("N34").Value = -62503400
$ws.Range("H58").Value = 3081.8635
$ws.Range("I58").Value = 3397.0557
$ws.Range("J58").Value = 1663.5
$ws.Range("K58").Value = 3397.0557
$ws.Range("L58").Value = 1663.5
$ws.Range("M58").Value = -3194.0557
$ws.Range("N58").Value = -2069.5
$ws.Range("H113").Value = 994.4167
$ws.Range("I113").Value = 1028.3334
$ws.Range("J113").Value = 892.6667
$ws.Range("K113").Value = 1028.3334
$ws.Range("L113").Value = 892.6667
$ws.Range("M113").Value = 1141.6666
$ws.Range("N113").Value = -5232.6667
$ws.Range("H136").Value = 3081.8635
$ws.Range("I136").Value = 3397.0557
$ws.Range("J136").Value = 1663.5
$ws.Range("K136").Value = 10191.1671
$ws.Range("L136").Value = 4990.5
$ws.Range("M136").Value = -7641.167099999999
$ws.Range("N136").Value = -10090.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1308.1666
$ws.Range("I7").Value = 169.8
$ws.Range("J7").Value = 7000
$ws.Range("K7").Value = 509.4
$ws.Range("L7").Value = 21000
$ws.Range("M7").Value = -397.4
$ws.Range("N7").Value = -21224
$ws.Range("H56").Value = 15977.5
$ws.Range("I56").Value = 15977.5
$ws.Range("K56").Value = 15977.5
$ws.Range("M56").Value = -15447.5
$ws.Range("H68").Value = 5643.4873
$ws.Range("I68").Value = 9999
$ws.Range("J68").Value = 5280.528
$ws.Range("K68").Value = 29997
$ws.Range("L68").Value = 15841.584
$ws.Range("M68").Value = -29186
$ws.Range("N68").Value = -17463.584
$ws.Range("H71").Value = 5643.4873
$ws.Range("I71").Value = 9999
$ws.Range("J71").Value = 5280.528
$ws.Range("K71").Value = 89991
$ws.Range("L71").Value = 47524.752
$ws.Range("M71").Value = -85935
$ws.Range("N71").Value = -55636.752
$ws.Range("H131").Value = 18524312
$ws.Range("I131").Value = 1196.875
$ws.Range("J131").Value = 33342804
$ws.Range("K131").Value = 3590.625
$ws.Range("L131").Value = 100028412
$ws.Range("M131").Value = 1449.375
$ws.Range("N131").Value = -100038492
$ws.Range("H138").Value = 2000
$ws.Range("I138").Value = 1000
$ws.Range("K138").Value = 3000
$ws.Range("M138").Value = 2140

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4560.5293
$ws.Range("I132").Value = 3902.8928
$ws.Range("J132").Value = 7629.5
$ws.Range("K132").Value = 11708.6784
$ws.Range("L132").Value = 22888.5
$ws.Range("M132").Value = -9178.678400000001
$ws.Range("N132").Value = -27948.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3474.2727
$ws.Range("J7").Value = 2998.5
$ws.Range("L7").Value = 2998.5
$ws.Range("N7").Value = -3222.5
$ws.Range("H22").Value = 1782.6923
$ws.Range("I22").Value = 1395.7778
$ws.Range("J22").Value = 2653.25
$ws.Range("K22").Value = 1395.7778
$ws.Range("L22").Value = 2653.25
$ws.Range("M22").Value = -1100.7778
$ws.Range("N22").Value = -3243.25
$ws.Range("H27").Value = 1782.6923
$ws.Range("I27").Value = 1395.7778
$ws.Range("J27").Value = 2653.25
$ws.Range("K27").Value = 1395.7778
$ws.Range("L27").Value = 2653.25
$ws.Range("M27").Value = -1288.7778
$ws.Range("N27").Value = -2867.25
$ws.Range("H126").Value = 3474.2727
$ws.Range("J126").Value = 2998.5
$ws.Range("L126").Value = 8995.5
$ws.Range("N126").Value = -13935.5
$ws.Range("H132").Value = 133335170
$ws.Range("I132").Value = 1860.5
$ws.Range("J132").Value = 400001800
$ws.Range("K132").Value = 5581.5
$ws.Range("L132").Value = 1200005400
$ws.Range("M132").Value = -3051.5
$ws.Range("N132").Value = -1200010460

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5528.5713
$ws.Range("I96").Value = 4950
$ws.Range("K96").Value = 4950
$ws.Range("M96").Value = -3577
$ws.Range("H126").Value = 7165.923
$ws.Range("I126").Value = 6217.4
$ws.Range("J126").Value = 10327.667
$ws.Range("K126").Value = 18652.2
$ws.Range("L126").Value = 30983.001
$ws.Range("M126").Value = -16182.2
$ws.Range("N126").Value = -35923.001
$ws.Range("H132").Value = 3398.6829
$ws.Range("I132").Value = 3693.147
$ws.Range("J132").Value = 1968.4286
$ws.Range("K132").Value = 11079.441
$ws.Range("L132").Value = 5905.2858
$ws.Range("M132").Value = -8549.440999999999
$ws.Range("N132").Value = -10965.2858
$ws.Range("H136").Value = 1350.2
$ws.Range("I136").Value = 1444.4375
$ws.Range("J136").Value = 1118.2307
$ws.Range("K136").Value = 4333.3125
$ws.Range("L136").Value = 3354.6921
$ws.Range("M136").Value = -1783.3125
$ws.Range("N136").Value = -8454.6921
